# The commit updates the big "score" callout on slide 1 from 0.678 to 0.66.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "0.678") {
            $shp.TextFrame.TextRange.Text = "0.66"
        }
    }
}
